$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores these as plain strings (prices such as
# "322.95", percentages padded with spaces such as "  -1.44%  ").
# Excel auto-converts number-looking text typed/assigned into a cell,
# so each literal below starts with a leading apostrophe (') to force
# text entry -- exactly like typing `'322.95` into the cell -- and the
# cell style is reset to "Normal" afterwards so the quote-prefix flag
# Excel records for that doesn't stick around (the source cells carry
# no explicit style).

$ws.Range("D2").Value = '''49.116.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.44%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.615.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.30%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.17%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''111.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.61%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''322.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -1.15%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -1.43%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.05%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.541'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -3.47%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''39.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -1.59%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''19.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -4.69%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.0809'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.24%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +1.06%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''7.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -0.34%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.030.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.03%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''2.622.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.23%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.855'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.96%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''49.099.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -1.32%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''3.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.11%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''12.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -3.46%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  -2.17%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.0₃0942'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.42%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''269.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -3.26%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''68.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -5.64%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -2.42%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''26.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.49%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.06%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  +3.50%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -0.30%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.138'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.93%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''34.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -4.93%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''49.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -0.39%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +0.29%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +2.35%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.11%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''18.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -4.01%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''4.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +3.18%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -1.30%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  +1.12%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''127.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.85%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -1.72%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''22.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -1.78%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.0320'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.35%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = '''2.055.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +0.21%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +7.10%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -4.39%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -9.41%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''8.89'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -1.43%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''58.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +1.82%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''5.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -3.20%  '
$ws.Range("E51").Style = "Normal"
